$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Section_A sheet
# ---------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

# Drop the now-unused trailing time slots (rows 9-12) first so the
# sheet dimension shrinks to A1:F8 once the remaining rows are filled in.
$wsA.Range("A9:F12").EntireRow.Delete()

$wsA.Range("A2").Value = "09:00-10:30"
$wsA.Range("B2").Value = "Free"
$wsA.Range("C2").Value = "MA102"
$wsA.Range("D2").Value = "MA102"
$wsA.Range("E2").Value = "HS101"
$wsA.Range("F2").Value = "EC101"

$wsA.Range("B3").Value = "HS101"
$wsA.Range("C3").Value = "DS101"
$wsA.Range("D3").Value = "EC101"
$wsA.Range("E3").Value = "Free"
$wsA.Range("F3").Value = "CS151 (Elective)"

$wsA.Range("A4").Value = "12:00-13:00"

$wsA.Range("A5").Value = "13:00-14:30"
$wsA.Range("B5").Value = "EC101"
$wsA.Range("C5").Value = "MA101"
$wsA.Range("D5").Value = "MA101"
$wsA.Range("E5").Value = "CS101"
$wsA.Range("F5").Value = "HS101"

$wsA.Range("A6").Value = "14:30-15:30"
$wsA.Range("B6").Value = "Free"
$wsA.Range("C6").Value = "Free"
$wsA.Range("E6").Value = "Free"

$wsA.Range("A7").Value = "15:30-17:00"
$wsA.Range("B7").Value = "Free"
$wsA.Range("C7").Value = "CS101"
$wsA.Range("D7").Value = "CS101"
$wsA.Range("E7").Value = "CS151 (Elective)"
$wsA.Range("F7").Value = "DS101"

$wsA.Range("A8").Value = "17:00-18:00"
$wsA.Range("C8").Value = "CS151 (Tutorial)"

# ---------------------------------------------------------------
# Section_B sheet
# ---------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("A9:F12").EntireRow.Delete()

$wsB.Range("A2").Value = "09:00-10:30"
$wsB.Range("B2").Value = "HS101"
$wsB.Range("C2").Value = "EC101"
$wsB.Range("D2").Value = "CS101"
$wsB.Range("E2").Value = "DS101"
$wsB.Range("F2").Value = "CS101"

$wsB.Range("D3").Value = "DS101"
$wsB.Range("E3").Value = "MA101"
$wsB.Range("F3").Value = "CS151 (Elective)"

$wsB.Range("A4").Value = "12:00-13:00"

$wsB.Range("A5").Value = "13:00-14:30"
$wsB.Range("B5").Value = "MA102"
$wsB.Range("C5").Value = "Free"
$wsB.Range("D5").Value = "EC101"
$wsB.Range("E5").Value = "CS101"
$wsB.Range("F5").Value = "MA101"

$wsB.Range("A6").Value = "14:30-15:30"
$wsB.Range("B6").Value = "Free"
$wsB.Range("C6").Value = "Free"
$wsB.Range("D6").Value = "Free"
$wsB.Range("E6").Value = "Free"
$wsB.Range("F6").Value = "Free"

$wsB.Range("A7").Value = "15:30-17:00"
$wsB.Range("B7").Value = "EC101"
$wsB.Range("C7").Value = "HS101"
$wsB.Range("D7").Value = "HS101"
$wsB.Range("E7").Value = "CS151 (Elective)"
$wsB.Range("F7").Value = "MA102"

$wsB.Range("A8").Value = "17:00-18:00"
$wsB.Range("C8").Value = "CS151 (Tutorial)"

# ---------------------------------------------------------------
# Course_Summary sheet - instructor reassignments
# ---------------------------------------------------------------
$wsC = $wb.Worksheets.Item("Course_Summary")

$wsC.Range("H2").Value = "Dr. Ramesh Adve"
$wsC.Range("H3").Value = "Dr. Abdul Wahid"
$wsC.Range("H5").Value = "Dr. Prakash Pawar"
$wsC.Range("H6").Value = "Dr. Sunil P V"
$wsC.Range("H8").Value = "Dr. Girish"

# ---------------------------------------------------------------
# Elective_Coordination sheet - expand with Session Type / Duration
# columns and break the elective into Lecture 1 / Lecture 2 / Tutorial
# ---------------------------------------------------------------
$wsE = $wb.Worksheets.Item("Elective_Coordination")

# F1 is a brand-new header cell (sheet used to stop at column E) -- give it
# the same header styling (bold, boxed, centered) as the rest of row 1
# before writing its text.
$wsE.Range("A1").Copy()
$wsE.Range("F1").PasteSpecial(-4122)
$wsE.Range("F1").Value = "Sections"

$wsE.Range("E1").Value = "Duration"
$wsE.Range("D1").Value = "Time Slot"
$wsE.Range("C1").Value = "Day"
$wsE.Range("B1").Value = "Session Type"

$wsE.Range("A2").Value = "CS151"
$wsE.Range("B2").Value = "Lecture 1"
$wsE.Range("C2").Value = "Fri"
$wsE.Range("D2").Value = "10:30-12:00"
$wsE.Range("E2").Value = "1.5 hours"
$wsE.Range("F2").Value = "A & B (Common Slot)"

$wsE.Range("A3").Value = "CS151"
$wsE.Range("B3").Value = "Lecture 2"
$wsE.Range("C3").Value = "Thu"
$wsE.Range("D3").Value = "15:30-17:00"
$wsE.Range("E3").Value = "1.5 hours"
$wsE.Range("F3").Value = "A & B (Common Slot)"

$wsE.Range("A4").Value = "CS151"
$wsE.Range("B4").Value = "Tutorial"
$wsE.Range("C4").Value = "Tue"
$wsE.Range("D4").Value = "17:00-18:00"
$wsE.Range("E4").Value = "1 hour"
$wsE.Range("F4").Value = "A & B (Common Slot)"
